# edit.ps1 - apply the resume edits described by the target diff.
#
# Strategy: use $d.Content.Find.Execute(...) to locate the anchor text for
# each change (never hard-coded character offsets, since those shift as
# earlier edits are applied), then use Range.InsertXML(...) to splice in
# the replacement OOXML with exact run/paragraph boundaries (Range.Text /
# Range.InsertBefore merges adjacent same-formatted runs together, which
# loses the run splits the target XML has).

$d = $word.ActiveDocument

function New-WordPackageXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        $bodyXml + `
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Insert a new "Current Coursework: ..." paragraph right after the
#    "Intro Signals & Systems..." / "...Intro Circuits" paragraph, ahead
#    of the existing lone <w:tab/> run that used to close that paragraph.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Intro Signals & Systems, Microprocessor Toys, Discrete Math, Intro Circuits") | Out-Null
$r.Collapse(0) | Out-Null              # wdCollapseEnd
$r.InsertParagraphAfter() | Out-Null   # splits off a new paragraph (same pPr), carrying the old tab run with it

$r.Collapse(0) | Out-Null
$r.Move(1, 1) | Out-Null               # step over the new paragraph mark, landing just before the tab run
$insertPos = $r.Start
$newParaRange = $d.Range($insertPos, $insertPos)
$bodyXml = '<w:p><w:r><w:t>Current Coursework: Operating Systems, Machine Learning (Coursera), Linear Algebra</w:t></w:r></w:p>'
$newParaRange.InsertXML((New-WordPackageXml $bodyXml)) | Out-Null

# ---------------------------------------------------------------------
# 2) Skills line: ", Git, LTSPICE" -> "," + _GoBack bookmark + " LTSPICE"
#    (Word keeps a single "_GoBack" bookmark marking the last edit spot,
#    so the one that used to sit after "setup " needs to move here.)
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete() | Out-Null
}

$r = $d.Content
$r.Find.Execute(", Git, LTSPICE") | Out-Null
$skillsRange = $d.Range($r.Start, $r.End)
$bodyXml = '<w:p><w:r><w:t>,</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t xml:space="preserve"> LTSPICE</w:t></w:r></w:p>'
$skillsRange.InsertXML((New-WordPackageXml $bodyXml)) | Out-Null

# ---------------------------------------------------------------------
# 3) Heading: "Project Experience" -> "P" + "ersonal Projects" (two runs)
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Project Experience") | Out-Null
$headingRange = $d.Range($r.Start, $r.End)
$bodyXml = '<w:p><w:r><w:t>P</w:t></w:r><w:r><w:t>ersonal Projects</w:t></w:r></w:p>'
$headingRange.InsertXML((New-WordPackageXml $bodyXml)) | Out-Null

# ---------------------------------------------------------------------
# 4) Merge the "Refactored signal processing...with" / "custom
#    gestures." paragraphs into a single paragraph + re-wrap as two runs.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Refactored signal processing algorithm to use machine learning to allow control of the device with") | Out-Null
$startPos = $r.Start

$r2 = $d.Content
$r2.Find.Execute("custom gestures.") | Out-Null
$endPos = $r2.End

$gestureRange = $d.Range($startPos, $endPos)
$bodyXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' + `
    '<w:tabs><w:tab w:val="right" w:pos="10800"/></w:tabs><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Refactored signal processing algorithm to use machine learning to </w:t></w:r>' + `
    '<w:r><w:t>recognize custom gestures.</w:t></w:r></w:p>'
$gestureRange.InsertXML((New-WordPackageXml $bodyXml)) | Out-Null

# ---------------------------------------------------------------------
# 5) Remove the two trailing Extracurriculars bullets: "International
#    Baccalaureate Diploma - Recipient" and "Michigan Club Wrestling -
#    Member", leaving "... Founder" as the last bullet.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("International Baccalaureate Diploma") | Out-Null
$startPos = $r.Start

$r2 = $d.Content
$r2.Find.Execute("Michigan Club Wrestling - Member") | Out-Null
$lastPara = $r2.Paragraphs(1)
$endPos = $lastPara.Range.End

$d.Range($startPos, $endPos).Delete() | Out-Null

Write-Output "edits applied"
